$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold header style, s="1") from column H into new columns I and J
$ws.Range("H1:H77").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1:I77").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row, I-value, J-value
$data = @(
    @(2, 6, 7),
    @(3, 7, 7),
    @(4, 6, 6),
    @(5, 7, 8),
    @(6, 9, 9),
    @(7, 10, 10),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 7, 7),
    @(12, 7, 7),
    @(13, 8, 8),
    @(14, 7, 8),
    @(15, 8, 9),
    @(16, 8, 8),
    @(17, 7, 7),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 7, 8),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 8, 8),
    @(29, 7, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 8, 8),
    @(36, 7, 7),
    @(37, 9, 9),
    @(38, 10, 10),
    @(39, 7, 7),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 7, 8),
    @(44, 8, 8),
    @(45, 7, 7),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 7, 7),
    @(49, 9, 9),
    @(50, 8, 8),
    @(51, 7, 7),
    @(52, 8, 8),
    @(53, 9, 9),
    @(54, 9, 9),
    @(55, 7, 7),
    @(56, 9, 9),
    @(57, 7, 7),
    @(58, 7, 9),
    @(59, 8, 8),
    @(60, 7, 8),
    @(61, 7, 8),
    @(62, 7, 8),
    @(63, 7, 7),
    @(64, 8, 8),
    @(65, 7, 8),
    @(66, 7, 8),
    @(67, 8, 8),
    @(68, 9, 9),
    @(69, 8, 8),
    @(70, 9, 9),
    @(71, 7, 8),
    @(72, 9, 9),
    @(73, 6, 6),
    @(74, 7, 7),
    @(75, 5, 5),
    @(76, 6, 6),
    @(77, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
